$d = $word.ActiveDocument
$d.Content.Find.Execute("Present time", $false, $false, $false, $false, $false, $true, 1, $false, "May 2019", 2)
